$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.383.86'
$ws.Range("E2").Value = '  -0.05%  '

# Row 3
$ws.Range("D3").Value = '1.846.94'
$ws.Range("E3").Value = '  -0.11%  '

# Row 4
$ws.Range("E4").Value = '  -0.27%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.47'
$ws.Range("E5").Value = '  -0.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6294'
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("E7").Value = '  -0.21%  '

# Row 8
$ws.Range("E8").Value = '  -1.65%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2909'
$ws.Range("E9").Value = '  +0.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.43'
$ws.Range("E10").Value = '  -1.15%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07712'
$ws.Range("E11").Value = '  -0.39%  '

# Row 12
$ws.Range("D12").Value = '1.846.77'
$ws.Range("E12").Value = '  -2.20%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.007'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6808'
$ws.Range("E14").Value = '  +0.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001029'
$ws.Range("E15").Value = '  -3.05%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.22'

# Row 17
$ws.Range("D17").Value = '2.102.54'
$ws.Range("E17").Value = '  -3.86%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.169'
$ws.Range("E18").Value = '  +0.31%  '

# Row 19
$ws.Range("D19").Value = '29.400.60'
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.22'
$ws.Range("E20").Value = '  +1.11%  '

# Row 21
$ws.Range("E21").Value = '  +0.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.457'
$ws.Range("E23").Value = '  -0.34%  '

# Row 24
$ws.Range("E24").Value = '  -0.23%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.74'
$ws.Range("E25").Value = '  +0.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1380'
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.419'
$ws.Range("E27").Value = '  +0.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.56'
$ws.Range("E28").Value = '  -0.62%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06450'
$ws.Range("E29").Value = '  +15.44%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.388'
$ws.Range("E30").Value = '  -0.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.475'
$ws.Range("E31").Value = '  +0.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.095'
$ws.Range("E32").Value = '  -0.69%  '

# Row 33
$ws.Range("E33").Value = '  +0.10%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.820'
$ws.Range("E34").Value = '  -0.81%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("E35").Value = '  -1.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6996'
$ws.Range("E36").Value = '  +0.27%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.577'
$ws.Range("E37").Value = '  -0.39%  '

# Row 38
$ws.Range("D38").Value = '1.261.56'
$ws.Range("E38").Value = '  +2.59%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.833'
$ws.Range("E39").Value = '  +4.19%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01822'
$ws.Range("E40").Value = '  +1.11%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.584'
$ws.Range("E41").Value = '  +2.80%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9091'
$ws.Range("E42").Value = '  +0.69%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9984'
$ws.Range("E43").Value = '  -0.27%  '

# Row 44
$ws.Range("D44").Value = '2.007.45'
$ws.Range("E44").Value = '  -18.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.41'
$ws.Range("E45").Value = '  -0.22%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.26'
$ws.Range("E46").Value = '  +0.55%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000119'
$ws.Range("E47").Value = '  +0.08%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1178'
$ws.Range("E48").Value = '  +2.88%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.065'
$ws.Range("E49").Value = '  -1.42%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.712'
$ws.Range("E50").Value = '  +2.20%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.059'
$ws.Range("E51").Value = '  +0.42%  '
